$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9 updates
$ws.Range("G9").Value = 3.4
$ws.Range("I9").Value = 2.25
$ws.Range("L9").Value = 3
$ws.Range("X9").Value = 15
$ws.Range("AS9").Value = 201
$ws.Range("AU9").Value = 8
$ws.Range("AW9").Value = 4.33

# Row 11 updates
$ws.Range("N11").Value = 8
